$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1 with the same style as H1 (bold/centered/bordered header)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-27
$values = @(
    @(6, 7),
    @(9, 9),
    @(6, 7),
    @(4, 6),
    @(7, 8),
    @(8, 8),
    @(7, 8),
    @(8, 9),
    @(7, 7),
    @(6, 7),
    @(6, 7),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(7, 8),
    @(6, 6),
    @(9, 9),
    @(7, 8),
    @(7, 8),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(2, 3),
    @(3, 3),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
